# Sort the curvature data (rows 2-12, columns A-D) by column A (time)
# ascending, leaving the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A2:D12")
$keyRange = $ws.Range("A2:A12")

# Header = xlNo (2) so the Sort call does not try to re-detect/skip a
# header row inside this already-headerless sub-range.
$dataRange.Sort($keyRange, 1, $null, $null, 1, $null, 1, 2)
